$d = $word.ActiveDocument

# --- Step 1: remove the existing "_GoBack" bookmark from the end of the
#     "... gemaakt" paragraph (paragraph 6). It will be re-created at the
#     end of the new final paragraph below.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- Step 2: append the new paragraphs after the trailing empty paragraph
#     (paragraph 7), i.e. right before the sectPr:
#       (empty)
#       24-1-2018
#       Contact pagina gemaakt in het eerste uur, ook het contactformulier
$lastEmpty = $d.Paragraphs.Last
$lastEmpty.Range.InsertParagraphAfter()

$dateParaIndex = $d.Paragraphs.Count
$dateParagraph = $d.Paragraphs.Item($dateParaIndex)
$dateParagraph.Range.Text = "24-1-2018"

$dateParagraph = $d.Paragraphs.Item($dateParaIndex)
$dateParagraph.Range.InsertParagraphAfter()

$contactParaIndex = $d.Paragraphs.Count
$contactParagraph = $d.Paragraphs.Item($contactParaIndex)
$contactParagraph.Range.Text = "Contact pagina gemaakt in het eerste uur, ook het contactformulier"

# --- Step 3: re-create the "_GoBack" bookmark at the end of the new last
#     paragraph's text (right before its paragraph mark), the same spot it
#     occupied relative to its original paragraph.
#     Placing a zero-width bookmark exactly at a paragraph's end position
#     is unreliable, so a throwaway marker character is appended first,
#     the bookmark is anchored just before it (a safe, non-boundary
#     position), and the marker character is then removed again.
$contactParagraph = $d.Paragraphs.Item($contactParaIndex)
$endPos = $contactParagraph.Range.End - 1
$markerRange = $d.Range($endPos, $endPos)
$markerRange.InsertAfter("#")

$contactParagraph = $d.Paragraphs.Item($contactParaIndex)
$bookmarkPos = $contactParagraph.Range.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$contactParagraph = $d.Paragraphs.Item($contactParaIndex)
$markerEnd = $contactParagraph.Range.End - 1
$removeRange = $d.Range($markerEnd - 1, $markerEnd)
$removeRange.Delete()
